# Increase the indentation of the Balance Sheet Items labels in column B
# so the hierarchy reads more clearly:
#   - Section headers / totals (currently 2 leading spaces) -> 4 leading spaces
#   - Line items (currently 4 leading spaces) -> 8 leading spaces
#   - Ratio labels (currently no leading spaces) -> 4 leading spaces

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section headers / subtotal rows: "  " -> "    "
$sectionRows = @(3, 9, 10, 15, 16, 19, 24, 31, 32, 36, 37, 41)
foreach ($r in $sectionRows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = "    " + $cell.Value2.TrimStart()
}

# Line-item rows: "    " -> "        "
$lineItemRows = @(4, 5, 6, 7, 8, 11, 12, 13, 14, 17, 18, 25, 26, 27, 28, 29, 30, 33, 34, 35, 38, 39, 40)
foreach ($r in $lineItemRows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = "        " + $cell.Value2.TrimStart()
}

# Ratio label rows (46-50): no leading spaces -> "    "
$ratioRows = @(46, 47, 48, 49, 50)
foreach ($r in $ratioRows) {
    $cell = $ws.Cells.Item($r, 2)
    $cell.Value2 = "    " + $cell.Value2.TrimStart()
}
